$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of "h1n1pdm" / "np" position data, appended after the existing
# data (which currently ends at row 133).
$positions = @(43, 87, 131, 175, 219, 263)

$startRow = 134
for ($i = 0; $i -lt $positions.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = "h1n1pdm"
    $ws.Cells.Item($r, 2).Value = "np"
    $ws.Cells.Item($r, 3).Value = $positions[$i]
}

$endRow = $startRow + $positions.Length - 1

# Match the updated selection recorded in the saved workbook.
$ws.Range("A${startRow}:C${endRow}").Select()
